# Updated cryptos list on Sun May 14 13:28:59 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.701.72"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.861.74"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("D4").Value = "'1.037"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.42%  "
$ws.Range("D5").Value = "'322.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").Value = "'1.033"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("D7").Value = "'0.4416"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("D8").Value = "'0.3790"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.07457"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("D10").Value = "'0.8829"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").Value = "'21.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("D12").Value = "1.873.84"
$ws.Range("E12").Value = "  -10.66%  "
$ws.Range("D13").Value = "'5.540"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").Value = "'6.753"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "'0.07217"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "'84.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("D17").Value = "'1.038"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "'0.000009107"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "27.704.29"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").Value = "'5.301"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").Value = "'11.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").Value = "2.089.77"
$ws.Range("E24").Value = "  -9.24%  "
$ws.Range("D25").Value = "'2.007"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.70%  "
$ws.Range("D26").Value = "'158.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").Value = "'18.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").Value = "'1.987"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("D29").Value = "'5.313"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("D30").Value = "'117.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").Value = "'0.7765"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.29%  "
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").Value = "'3.019"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.35%  "
$ws.Range("D35").Value = "'4.570"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.49%  "
$ws.Range("D36").Value = "'1.034"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("D39").Value = "'0.05336"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("D40").Value = "'2.866"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.83%  "
$ws.Range("D41").Value = "'0.5191"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("D43").Value = "'6.844"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.88%  "
$ws.Range("D44").Value = "'8.657"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.49%  "
$ws.Range("D45").Value = "'110.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("D46").Value = "'10.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").Value = "'0.06633"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.72%  "
$ws.Range("D48").Value = "'1.711"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.95%  "
$ws.Range("D49").Value = "'0.4706"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").Value = "'1.910"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").Value = "'39.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.46%  "
